$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Update the "No" column (A) with new sequential values starting at 33
$noValues = 33..42
for ($i = 0; $i -lt $noValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $noValues[$i]
}

# Update "Created At" (E) and "Update At" (F) columns for rows 2-11 to the new timestamp
$newTimestamp = "2021-05-03T14:53:17.000000Z"
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 5).Value = $newTimestamp
    $ws.Cells.Item($row, 6).Value = $newTimestamp
}
